# Making the Case for Rust - apply author's latest round of edits
#
# 1. Insert a new "Collections" slide (Title and Content layout) right
#    before the "For C Developers" slide (currently slide index 13).
# 2. Merge the two title runs on the "Levels of IoT Devices" slide
#    (slide 2) into a single run.

$p = $ppt.ActivePresentation

# --- 1. New "Collections" slide -------------------------------------------
# "Title and Content" is CustomLayout #2 on the slide master (same layout
# used by the neighbouring slides such as "For C Developers").
$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.Add(13, 2)
$newSlide.Design = $p.SlideMaster.Design

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Collections"

# --- 2. Tidy up the "Levels of IoT Devices" title on slide 2 --------------
$s2 = $p.Slides.Item(2)
$titleShape = $s2.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
# Force a real text replacement (the concatenated text is already correct,
# so round-trip through an unrelated placeholder value first) so the two
# separate runs collapse into a single run.
$titleRange.Text = "x"
$titleRange.Text = "Levels of IoT Devices"
